# Applies the commit's changes:
#  1. Rename "Weekly Quantity"!B1  "Requested quantity" -> "Weekly_PO_Qty"
#  2. Rename "Monthly Trend"!B1    "Requested quantity" -> "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing a
#     ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

# --- 1. Update "Weekly Quantity" sheet header ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Update "Monthly Trend" sheet header ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page layout / outline defaults used by the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Copy header formatting (bold, centered, bordered) from the weekly sheet's header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column formatting (custom date/time number format) down column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A45").PasteSpecial(-4122)

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Bulk-write the forecast data (44 rows x 4 cols) into A2:D45
$data = New-Object 'object[,]' 44,4
$data[0,0] = 45011.99999999999
$data[0,1] = 4
$data[0,2] = -1.157394087867994
$data[0,3] = 10.17622008563799
$data[1,0] = 45025.99999999999
$data[1,1] = 4
$data[1,2] = -0.7084387095703976
$data[1,3] = 10.32939685484788
$data[2,0] = 45032.99999999999
$data[2,1] = 4
$data[2,2] = -0.8767252070689162
$data[2,3] = 9.665739274752104
$data[3,0] = 45039.99999999999
$data[3,1] = 4
$data[3,2] = -0.9983434485576641
$data[3,3] = 9.745191679083058
$data[4,0] = 45046.99999999999
$data[4,1] = 4
$data[4,2] = -1.453431349080611
$data[4,3] = 10.06101943356497
$data[5,0] = 45053.99999999999
$data[5,1] = 4
$data[5,2] = -0.8171259912851397
$data[5,3] = 10.01519349822399
$data[6,0] = 45060.99999999999
$data[6,1] = 5
$data[6,2] = -0.9789737211462721
$data[6,3] = 9.905053700305457
$data[7,0] = 45123.99999999999
$data[7,1] = 5
$data[7,2] = -0.8562099680459888
$data[7,3] = 10.50337642528063
$data[8,0] = 45130.99999999999
$data[8,1] = 5
$data[8,2] = -1.240394351750074
$data[8,3] = 10.26486005762117
$data[9,0] = 45137.99999999999
$data[9,1] = 5
$data[9,2] = -0.7294744196614545
$data[9,3] = 10.53704002771888
$data[10,0] = 45144.99999999999
$data[10,1] = 5
$data[10,2] = -0.7629665479839044
$data[10,3] = 10.55693927992275
$data[11,0] = 45151.99999999999
$data[11,1] = 5
$data[11,2] = -0.9382557512362532
$data[11,3] = 10.05226201763131
$data[12,0] = 45312.99999999999
$data[12,1] = 5
$data[12,2] = -0.2981270141250365
$data[12,3] = 10.71153163887824
$data[13,0] = 45326.99999999999
$data[13,1] = 5
$data[13,2] = -0.4211458230581507
$data[13,3] = 10.8904078515062
$data[14,0] = 45333.99999999999
$data[14,1] = 5
$data[14,2] = -0.1135917635970542
$data[14,3] = 10.47284491957235
$data[15,0] = 45347.99999999999
$data[15,1] = 5
$data[15,2] = -0.2383968226852416
$data[15,3] = 11.20076141853068
$data[16,0] = 45354.99999999999
$data[16,1] = 6
$data[16,2] = 0.1586697304639554
$data[16,3] = 11.14065189682535
$data[17,0] = 45361.99999999999
$data[17,1] = 6
$data[17,2] = 0.6050317287650741
$data[17,3] = 11.08819660373119
$data[18,0] = 45368.99999999999
$data[18,1] = 6
$data[18,2] = -0.1838027639288522
$data[18,3] = 10.6840396648521
$data[19,0] = 45375.99999999999
$data[19,1] = 6
$data[19,2] = 0.03663544076195988
$data[19,3] = 11.18564973539178
$data[20,0] = 45382.99999999999
$data[20,1] = 6
$data[20,2] = -0.05364784477833526
$data[20,3] = 10.97110904817463
$data[21,0] = 45389.99999999999
$data[21,1] = 6
$data[21,2] = 0.5015113129917435
$data[21,3] = 11.0534007901653
$data[22,0] = 45396.99999999999
$data[22,1] = 6
$data[22,2] = -0.1690308091786527
$data[22,3] = 11.00640460999188
$data[23,0] = 45459.99999999999
$data[23,1] = 6
$data[23,2] = 0.5897612912377753
$data[23,3] = 11.36888760289435
$data[24,0] = 45487.99999999999
$data[24,1] = 6
$data[24,2] = 0.8171149282170486
$data[24,3] = 11.34837260302089
$data[25,0] = 45515.99999999999
$data[25,1] = 6
$data[25,2] = 0.8621975167393384
$data[25,3] = 11.51430980497509
$data[26,0] = 45522.99999999999
$data[26,1] = 6
$data[26,2] = 0.3458742974674109
$data[26,3] = 11.16672655770474
$data[27,0] = 45536.99999999999
$data[27,1] = 6
$data[27,2] = 0.6878288220320558
$data[27,3] = 11.38475982040173
$data[28,0] = 45543.99999999999
$data[28,1] = 6
$data[28,2] = 0.9371934499598081
$data[28,3] = 11.73733104447085
$data[29,0] = 45550.99999999999
$data[29,1] = 6
$data[29,2] = 0.3681957494194085
$data[29,3] = 11.80531112625285
$data[30,0] = 45557.99999999999
$data[30,1] = 6
$data[30,2] = 1.201348783045256
$data[30,3] = 11.88801364617314
$data[31,0] = 45564.99999999999
$data[31,1] = 6
$data[31,2] = 0.8930776438032432
$data[31,3] = 11.64855804988502
$data[32,0] = 45571.99999999999
$data[32,1] = 6
$data[32,2] = 0.6332314142957464
$data[32,3] = 11.77065908395544
$data[33,0] = 45578.99999999999
$data[33,1] = 6
$data[33,2] = 1.307165857051223
$data[33,3] = 11.7035209703224
$data[34,0] = 45585.99999999999
$data[34,1] = 6
$data[34,2] = 1.165827034566659
$data[34,3] = 11.84415811490669
$data[35,0] = 45592.99999999999
$data[35,1] = 6
$data[35,2] = 0.5563662151954832
$data[35,3] = 11.83735197270789
$data[36,0] = 45599.99999999999
$data[36,1] = 6
$data[36,2] = 0.8964408811613213
$data[36,3] = 12.0834144649038
$data[37,0] = 45606.99999999999
$data[37,1] = 6
$data[37,2] = 0.8952131400953894
$data[37,3] = 11.7651366467835
$data[38,0] = 45613.99999999999
$data[38,1] = 6
$data[38,2] = 1.241739306612066
$data[38,3] = 11.70188438351058
$data[39,0] = 45620.99999999999
$data[39,1] = 6
$data[39,2] = 1.073029688754216
$data[39,3] = 12.01603233260422
$data[40,0] = 45627.99999999999
$data[40,1] = 6
$data[40,2] = 1.025340110207794
$data[40,3] = 11.82045979992665
$data[41,0] = 45634.99999999999
$data[41,1] = 6
$data[41,2] = 1.197457736639431
$data[41,3] = 11.93018871100849
$data[42,0] = 45641.99999999999
$data[42,1] = 7
$data[42,2] = 1.261857042238575
$data[42,3] = 12.12111417584256
$data[43,0] = 45648.99999999999
$data[43,1] = 7
$data[43,2] = 0.9802135886251143
$data[43,3] = 11.9894250635266

$wsForecast.Range("A2:D45").Value = $data

# Restore the originally active sheet/tab
$wsWeekly.Activate()
